$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain number-looking string need NumberFormat
# forced to Text ("@") first, otherwise COM auto-converts the assigned
# string into a numeric cell and the literal text (e.g. trailing zeros) is lost.

$ws.Range("D2").Value = '64.021.09'
$ws.Range("E2").Value = '  -2.46%  '

$ws.Range("D3").Value = '3.476.91'
$ws.Range("E3").Value = '  -3.02%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.61'
$ws.Range("E5").Value = '  -3.37%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.13'
$ws.Range("E6").Value = '  -4.16%  '

$ws.Range("D7").Value = '3.472.72'
$ws.Range("E7").Value = '  -3.14%  '

$ws.Range("E8").Value = '  +0.14%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.488'
$ws.Range("E9").Value = '  -1.85%  '

$ws.Range("E10").Value = '  -1.34%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.18'
$ws.Range("E11").Value = '  -1.03%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.384'
$ws.Range("E12").Value = '  -1.92%  '

$ws.Range("D13").Value = '4.071.24'
$ws.Range("E13").Value = '  -3.02%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.69'
$ws.Range("E14").Value = '  -1.28%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("B15").Value = 'TRON'
$ws.Range("C15").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D15").Value = '0.117'
$ws.Range("E15").Value = '  +0.14%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '0.0000177'
$ws.Range("E16").Value = '  -5.07%  '

$ws.Range("D17").Value = '3.483.07'
$ws.Range("E17").Value = '  -2.96%  '

$ws.Range("D18").Value = '63.580.72'
$ws.Range("E18").Value = '  -1.91%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.90'
$ws.Range("E19").Value = '  -1.21%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.21'
$ws.Range("E20").Value = '  -2.75%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.63'
$ws.Range("E21").Value = '  -4.59%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '392.04'
$ws.Range("E22").Value = '  -1.30%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.573'
$ws.Range("E23").Value = '  -2.75%  '

$ws.Range("D24").Value = '3.619.30'
$ws.Range("E24").Value = '  -3.01%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '73.18'
$ws.Range("E25").Value = '  -1.48%  '

$ws.Range("E26").Value = '  +0.09%  '

$ws.Range("E27").Value = '  -8.67%  '

$ws.Range("E28").Value = '  -5.92%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.46'
$ws.Range("E29").Value = '  -8.67%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.989'
$ws.Range("E30").Value = '  -1.14%  '

$ws.Range("E31").Value = '  -7.11%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.14'
$ws.Range("E32").Value = '  -5.29%  '

$ws.Range("D33").Value = '3.476.33'
$ws.Range("E33").Value = '  -3.16%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.76'
$ws.Range("E35").Value = '  -3.08%  '

$ws.Range("E36").Value = '  -2.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.22'
$ws.Range("E37").Value = '  -3.35%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").Value = '169.67'
$ws.Range("E38").Value = '  +0.84%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").Value = '1.56'
$ws.Range("E39").Value = '  -2.64%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").Value = '6.94'
$ws.Range("E40").Value = '  -2.11%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0802'
$ws.Range("E41").Value = '  -4.24%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.809'
$ws.Range("E42").Value = '  -3.82%  '

$ws.Range("E43").Value = '  -3.78%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  +0.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.55'
$ws.Range("E45").Value = '  -3.94%  '

$ws.Range("E46").Value = '  -6.44%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.34'
$ws.Range("E47").Value = '  -4.52%  '

$ws.Range("E48").Value = '  -3.93%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.84'
$ws.Range("E49").Value = '  -2.83%  '

$ws.Range("D50").Value = '2.416.80'
$ws.Range("E50").Value = '  -1.30%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0265'
$ws.Range("E51").Value = '  -2.11%  '

